$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text labels: replace the space between words with a line break ---
$ws.Range("B1").Value = "Kia`nRio"
$ws.Range("C1").Value = "Volkswagen`nGolf"
$ws.Range("D1").Value = "Toyota`nCorolla"
$ws.Range("E1").Value = "Skoda`nOctavia"
$ws.Range("F1").Value = "BMW`n3`nSeries"
$ws.Range("G1").Value = "Hyundai`nSolaris"
$ws.Range("H1").Value = "Вектор`nприоритетов"

$ws.Range("A2").Value = "Kia`nRio"
$ws.Range("A3").Value = "Volkswagen`nGolf"
$ws.Range("A4").Value = "Toyota`nCorolla"
$ws.Range("A5").Value = "Skoda`nOctavia"
$ws.Range("A6").Value = "BMW`n3`nSeries"
$ws.Range("A7").Value = "Hyundai`nSolaris"

# --- Update slightly recalculated priority vector values (kept as text) ---
$ws.Range("H3").Formula = "=""0.228"""
$ws.Range("H3").Copy()
$ws.Range("H3").PasteSpecial(-4163)
$ws.Range("H6").Formula = "=""0.362"""
$ws.Range("H6").Copy()
$ws.Range("H6").PasteSpecial(-4163)

# --- Update column widths ---
$ws.Columns.Item(1).ColumnWidth = 14.4
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(3).ColumnWidth = 14.4
$ws.Columns.Item(4).ColumnWidth = 10.8
$ws.Columns.Item(5).ColumnWidth = 10.8
$ws.Columns.Item(6).ColumnWidth = 9.6
$ws.Columns.Item(7).ColumnWidth = 10.8
$ws.Columns.Item(8).ColumnWidth = 15.6
